# Weekly update: insert 3 new rows of data (new week, date 44419) right
# before the existing "Verde" quality-grade rows for Vega Central
# Mapocho de Santiago - Espárragos, pushing the old rows 9-17 down to 12-20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at row 9 - this shifts existing rows 9:17 down to 12:20.
$ws.Rows("9:11").Insert()

# Row 9 - Banquete, $/bandeja 10 kilos
$ws.Range("A9").Value = 9
$ws.Range("B9").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C9").Value = "Metropolitana"
$ws.Range("D9").Value = 44419
$ws.Range("E9").Value = 13
$ws.Range("F9").Value = 300000000
$ws.Range("G9").Value = "Espárragos"
$ws.Range("H9").Value = "Verde"
$ws.Range("I9").Value = "Banquete"
$ws.Range("J9").Value = 7
$ws.Range("K9").Value = 35000
$ws.Range("L9").Value = 35000
$ws.Range("M9").Value = 35000
$ws.Range("N9").Value = "$/bandeja 10 kilos"
$ws.Range("O9").Value = "Región Metropolitana"
$ws.Range("P9").Value = 3500
$ws.Range("Q9").Value = 10
$ws.Range("R9").Value = "Hortaliza"

# Row 10 - Primera, $/bandeja 10 kilos
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C10").Value = "Metropolitana"
$ws.Range("D10").Value = 44419
$ws.Range("E10").Value = 13
$ws.Range("F10").Value = 300000000
$ws.Range("G10").Value = "Espárragos"
$ws.Range("H10").Value = "Verde"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 16
$ws.Range("K10").Value = 34000
$ws.Range("L10").Value = 34000
$ws.Range("M10").Value = 34000
$ws.Range("N10").Value = "$/bandeja 10 kilos"
$ws.Range("O10").Value = "Región Metropolitana"
$ws.Range("P10").Value = 3400
$ws.Range("Q10").Value = 10
$ws.Range("R10").Value = "Hortaliza"

# Row 11 - Segunda, $/bandeja 10 kilos
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C11").Value = "Metropolitana"
$ws.Range("D11").Value = 44419
$ws.Range("E11").Value = 13
$ws.Range("F11").Value = 300000000
$ws.Range("G11").Value = "Espárragos"
$ws.Range("H11").Value = "Verde"
$ws.Range("I11").Value = "Segunda"
$ws.Range("J11").Value = 9
$ws.Range("K11").Value = 32000
$ws.Range("L11").Value = 32000
$ws.Range("M11").Value = 32000
$ws.Range("N11").Value = "$/bandeja 10 kilos"
$ws.Range("O11").Value = "Región Metropolitana"
$ws.Range("P11").Value = 3200
$ws.Range("Q11").Value = 10
$ws.Range("R11").Value = "Hortaliza"
